$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 3.25
$ws.Range("I2").Value = 2.45
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 7.5
$ws.Range("Y2").Value = 1.53
$ws.Range("Z2").Value = 2.38
$ws.Range("AC2").Value = 8.5
$ws.Range("AI2").Value = 7.5
$ws.Range("AK2").Value = 15
$ws.Range("AN2").Value = 7

# Row 4 updates
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("AB4").Value = 1.69

# Row 5 updates
$ws.Range("G5").Value = 1.55
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.2
$ws.Range("K5").Value = 2.2
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("S5").Value = 2.05
$ws.Range("T5").Value = 1.75
$ws.Range("W5").Value = 3.75
$ws.Range("X5").Value = 1.25
$ws.Range("AB5").Value = 1.63
$ws.Range("AC5").Value = 6
$ws.Range("AD5").Value = 7
$ws.Range("AE5").Value = 8.5
$ws.Range("AH5").Value = 29
$ws.Range("AI5").Value = 9
$ws.Range("AJ5").Value = 7.5
$ws.Range("AK5").Value = 19
$ws.Range("AN5").Value = 13
